$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMSIN")

# A completely untouched / never-formatted cell used purely as a "blank
# format" donor so PasteSpecial(formats) can strip an unwanted auto-assigned
# style back down to "no override" on cells that must stay at the sheet's
# default format.
$blankDonor = $ws.Range("ZZ1")

# ---------------------------------------------------------------------------
# Row 58: existing run (2022-09-02 / s166) picks up the plain "data" style
# on every column except the already-styled datetime cell B58.
# ---------------------------------------------------------------------------

# B58 keeps its existing datetime style; only the stored value shifts by a
# hair of floating-point precision (same instant, re-serialised).
$ws.Range("B58").Value = 44806.6996547338

# C58/D58/E58/F58/G58: re-assert the same values through Clear+Value so the
# cells pick up the workbook's normal populated-cell style.
$ws.Range("C58").ClearContents()
$ws.Range("C58").Value = "s166"
$ws.Range("D58").ClearContents()
$ws.Range("D58").Value = 51
$ws.Range("E58").ClearContents()
$ws.Range("E58").Value = 51
$ws.Range("F58").ClearContents()
$ws.Range("F58").Value = 0
$ws.Range("G58").ClearContents()
$ws.Range("G58").Value = 1.07

# A58 holds a date-shaped string ("2022-09-02") that must stay literal text,
# not get silently converted into a serial date. Writing it with a leading
# apostrophe forces text, then we apply the same style used on the rest of
# the row so every cell in row 58 is formatted consistently.
$rowStyle = $ws.Range("D58").Style
$ws.Range("A58").Value = "'2022-09-02"
$ws.Range("A58").Style = $rowStyle

# ---------------------------------------------------------------------------
# Row 59: brand-new certificate run. Target cells keep the sheet's default
# (unstyled) look except B59, which adopts B58's datetime style.
# ---------------------------------------------------------------------------

# B59 already exists (was a blank date-styled placeholder) - copy B58's
# format across so it matches the same datetime style, then set its value.
$ws.Range("B58").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$ws.Range("B59").Value = 44810.9392168561

# C59/D59/E59/F59/G59 are new cells; a plain value write auto-promotes them
# to the populated-cell style, so strip that back down to the sheet default
# via a format-only paste from an untouched cell.
$ws.Range("C59").Value = "cert89"
$blankDonor.Copy()
$ws.Range("C59").PasteSpecial(-4122)

$ws.Range("D59").Value = 51
$blankDonor.Copy()
$ws.Range("D59").PasteSpecial(-4122)

$ws.Range("E59").Value = 51
$blankDonor.Copy()
$ws.Range("E59").PasteSpecial(-4122)

$ws.Range("F59").Value = 0
$blankDonor.Copy()
$ws.Range("F59").PasteSpecial(-4122)

$ws.Range("G59").Value = 1.15
$blankDonor.Copy()
$ws.Range("G59").PasteSpecial(-4122)

# A59 is the same "date-shaped text" situation as A58: force text with a
# leading apostrophe, then strip the resulting quote-prefix format back to
# the sheet default (row 59's new cells stay unstyled in the target).
$ws.Range("A59").Value = "'2022-09-06"
$blankDonor.Copy()
$ws.Range("A59").PasteSpecial(-4122)
